$p = $ppt.ActivePresentation

# --- Add a second slide using the same (Blank) layout as slide 1 ---
$layout = $p.SlideMaster.CustomLayouts.Item(1)
$s2 = $p.Slides.AddSlide(2, $layout)

# --- Slide 1: notes page with text in the inherited position/size ---
$s1 = $p.Slides.Item(1)
$notes1 = $s1.NotesPage
$notesBody1 = $notes1.Shapes.AddPlaceholder(2)
$notesBody1.TextFrame.TextRange.Text = "Notes in inherited position and size."

# --- Slide 2: notes page with text in a modified position/size ---
$notes2 = $s2.NotesPage
$notesBody2 = $notes2.Shapes.AddPlaceholder(2)
$notesBody2.TextFrame.TextRange.Text = "Notes in modified position and size."
